$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet stays untouched. Only the per-language report sheets
# (zh-cn / de-de) get the "handback" columns (Latest Target File / Latest
# Handback File) populated, plus the handoff-status text + handback
# timestamps get refreshed.
# ---------------------------------------------------------------------------

$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

# --- zh-cn sheet: fully handed back, in sync with en-US -------------------

# Row 2 (6e1e0346-... file)
$zh.Hyperlinks.Add(
    $zh.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/633b24567c09d380ffeeac20bff3e2d73abbaac6/e2e/6e1e0346-9eaf-4344-a994-a53c0294167c.md",
    [Type]::Missing,
    [Type]::Missing,
    "6e1e0346-9eaf-4344-a994-a53c0294167c.md"
)
$zh.Hyperlinks.Add(
    $zh.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6e64b103b9750747d40696cb9561844ae5b55365/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/6e1e0346-9eaf-4344-a994-a53c0294167c.0e6725bc6deaf82c7bfc95a4fc6fb1be945f61e4.zh-cn.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "6e1e0346-9eaf-4344-a994-a53c0294167c.0e6725bc6deaf82c7bfc95a4fc6fb1be945f61e4.zh-cn.xlf"
)
$zh.Range("H2").Value = "2016-03-14 09:43:08"

# Row 3 (f37793fb-... file)
$zh.Hyperlinks.Add(
    $zh.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/633b24567c09d380ffeeac20bff3e2d73abbaac6/e2e/f37793fb-955e-4c6f-ac70-5f46187ab8df.md",
    [Type]::Missing,
    [Type]::Missing,
    "f37793fb-955e-4c6f-ac70-5f46187ab8df.md"
)
$zh.Hyperlinks.Add(
    $zh.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6e64b103b9750747d40696cb9561844ae5b55365/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht/f37793fb-955e-4c6f-ac70-5f46187ab8df.fde37e87e2b3d3e56dc8b08a95d60d2afbd03a37.zh-cn.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "f37793fb-955e-4c6f-ac70-5f46187ab8df.fde37e87e2b3d3e56dc8b08a95d60d2afbd03a37.zh-cn.xlf"
)
$zh.Range("H3").Value = "2016-03-14 09:43:08"

# --- de-de sheet: also handed back, in sync with en-US --------------------

# Row 2 (6e1e0346-... file)
$de.Hyperlinks.Add(
    $de.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/633b24567c09d380ffeeac20bff3e2d73abbaac6/e2e/6e1e0346-9eaf-4344-a994-a53c0294167c.md",
    [Type]::Missing,
    [Type]::Missing,
    "6e1e0346-9eaf-4344-a994-a53c0294167c.md"
)
$de.Hyperlinks.Add(
    $de.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9446d685db8777377af652cbaca689e6699ec644/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/6e1e0346-9eaf-4344-a994-a53c0294167c.0e6725bc6deaf82c7bfc95a4fc6fb1be945f61e4.de-de.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "6e1e0346-9eaf-4344-a994-a53c0294167c.0e6725bc6deaf82c7bfc95a4fc6fb1be945f61e4.de-de.xlf"
)
$de.Range("H2").Value = "2016-03-14 09:43:21"

# Row 3 (f37793fb-... file)
$de.Hyperlinks.Add(
    $de.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/633b24567c09d380ffeeac20bff3e2d73abbaac6/e2e/f37793fb-955e-4c6f-ac70-5f46187ab8df.md",
    [Type]::Missing,
    [Type]::Missing,
    "f37793fb-955e-4c6f-ac70-5f46187ab8df.md"
)
$de.Hyperlinks.Add(
    $de.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9446d685db8777377af652cbaca689e6699ec644/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht/f37793fb-955e-4c6f-ac70-5f46187ab8df.fde37e87e2b3d3e56dc8b08a95d60d2afbd03a37.de-de.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "f37793fb-955e-4c6f-ac70-5f46187ab8df.fde37e87e2b3d3e56dc8b08a95d60d2afbd03a37.de-de.xlf"
)
$de.Range("H3").Value = "2016-03-14 09:43:21"

# --- Status text refresh ---------------------------------------------------
# "Ready for handoff" -> "Handed back: in sync with en-US" everywhere the
# status is shown: the Overview sheet's zh-cn/de-de columns (B,C) for both
# rows, and the per-language sheets' Status column (C) for both rows.

$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = "Handed back: in sync with en-US"
$ov.Range("C2").Value = "Handed back: in sync with en-US"
$ov.Range("B3").Value = "Handed back: in sync with en-US"
$ov.Range("C3").Value = "Handed back: in sync with en-US"

$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("C3").Value = "Handed back: in sync with en-US"

$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("C3").Value = "Handed back: in sync with en-US"
